# Rerun all TODE norms, found error in grade code that was creating the weird column
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("K-Spring")
$ws.Range("B8").Value = 79
$ws.Range("B24").Value = 124

$ws = $wb.Worksheets.Item("1-Fall")
$ws.Range("B2").Value = 52
$ws.Range("B3").Value = 56
$ws.Range("B26").Value = 122

$ws = $wb.Worksheets.Item("1-Spring")
$ws.Range("B2").Value = 42
$ws.Range("B3").Value = 46
$ws.Range("B4").Value = 50
$ws.Range("B5").Value = 54
$ws.Range("B6").Value = 58
$ws.Range("B7").Value = 61
$ws.Range("B8").Value = 64
$ws.Range("B9").Value = 67
$ws.Range("B10").Value = 70
$ws.Range("B11").Value = 73
$ws.Range("B12").Value = 76
$ws.Range("B13").Value = 79
$ws.Range("B14").Value = 82
$ws.Range("B15").Value = 84
$ws.Range("B16").Value = 87
$ws.Range("B17").Value = 90
$ws.Range("B18").Value = 93
$ws.Range("B19").Value = 95
$ws.Range("B20").Value = 98
$ws.Range("B21").Value = 101
$ws.Range("B22").Value = 103
$ws.Range("B23").Value = 106
$ws.Range("B24").Value = 109
$ws.Range("B25").Value = 112
$ws.Range("B26").Value = 115
$ws.Range("B27").Value = 118
$ws.Range("B28").Value = 121
$ws.Range("B29").Value = 124
$ws.Range("B30").Value = 128

$ws = $wb.Worksheets.Item("2-Fall")
$ws.Range("B4").Value = 40
$ws.Range("B5").Value = 44
$ws.Range("B6").Value = 48
$ws.Range("B7").Value = 52
$ws.Range("B8").Value = 55
$ws.Range("B9").Value = 59
$ws.Range("B10").Value = 62
$ws.Range("B11").Value = 65
$ws.Range("B22").Value = 97
$ws.Range("B27").Value = 110
$ws.Range("B28").Value = 113
$ws.Range("B31").Value = 123

$ws = $wb.Worksheets.Item("2-Spring")
$ws.Range("B7").Value = 41
$ws.Range("B8").Value = 45
$ws.Range("B9").Value = 49
$ws.Range("B10").Value = 53
$ws.Range("B11").Value = 57
$ws.Range("B12").Value = 60
